$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.00000396693397741777018937654741836951188815874047577381134
$ws.Range("E2").Value = 0.00000396693397741777018937654741836951188815874047577381134

$ws.Range("D3").Value = 0.999999999992604582388366907252930104732513427734375
$ws.Range("E3").Value = 0.999999999992604582388366907252930104732513427734375

$ws.Range("D4").Value = 0.99852484647533767958549333343398757278919219970703125
$ws.Range("E4").Value = 0.99852484647533767958549333343398757278919219970703125

$ws.Range("D5").Value = 0.000000000000000000000000000000000000000372911067838235494115
$ws.Range("E5").Value = 0.000000000000000000000000000000000000000372911067838235494115

$ws.Range("D6").Value = 0.007078098031594035229041583789921787683852016925811767578125
$ws.Range("E6").Value = 0.007078098031594035229041583789921787683852016925811767578125

$ws.Range("D7").Value = 0.9999999953767737537191351293586194515228271484375
$ws.Range("E7").Value = 0.0000000046232262462808648706413805484771728515625

$ws.Range("D8").Value = 0.000051197709181899523158982417969653511136129964143037796021
$ws.Range("E8").Value = 0.999948802290818061777599723427556455135345458984375

$ws.Range("D9").Value = 0.99026783246948124617148323522997088730335235595703125
$ws.Range("E9").Value = 0.00973216753051875382851676477002911269664764404296875

$ws.Range("D11").Value = 0.999999999476063550218896125443279743194580078125
$ws.Range("E11").Value = 0.000000000523936449781103874556720256805419921875
$ws.Range("F11").Value = 4.20458507537841796875
